# Update ra/dec residual columns (AF:AI) with newly computed variance-based values
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("AF2").Value = [double]"-5.389228604713026e-06"
$ws.Range("AG2").Value = [double]"-1.777760231735925e-05"
$ws.Range("AH2").Value = [double]"-9.405978329490228e-08"
$ws.Range("AI2").Value = [double]"-3.102776935480928e-07"

$ws.Range("AF3").Value = [double]"0.0001562398980325952"
$ws.Range("AG3").Value = [double]"0.0001237796038484618"
$ws.Range("AH3").Value = [double]"2.726900643648997e-06"
$ws.Range("AI3").Value = [double]"2.160361633969902e-06"

$ws.Range("AF4").Value = [double]"5.124683426060983e-06"
$ws.Range("AG4").Value = [double]"-3.252834425460094e-05"
$ws.Range("AH4").Value = [double]"8.944259890714754e-08"
$ws.Range("AI4").Value = [double]"-5.677267074649671e-07"

$ws.Range("AF5").Value = [double]"5.512698855625331e-05"
$ws.Range("AG5").Value = [double]"-0.0001353706841280555"
$ws.Range("AH5").Value = [double]"9.621474570158554e-07"
$ws.Range("AI5").Value = [double]"-2.362664148711798e-06"

$ws.Range("AF6").Value = [double]"4.34363420254158e-06"
$ws.Range("AG6").Value = [double]"-2.581056363304413e-06"
$ws.Range("AH6").Value = [double]"7.581071833658882e-08"
$ws.Range("AI6").Value = [double]"-4.504793171921295e-08"

$ws.Range("AF7").Value = [double]"4.630711065090054e-05"
$ws.Range("AG7").Value = [double]"7.675700004217134e-05"
$ws.Range("AH7").Value = [double]"8.082115479435489e-07"
$ws.Range("AI7").Value = [double]"1.339662374689316e-06"

$ws.Range("AF8").Value = [double]"-6.822983223742085e-06"
$ws.Range("AG8").Value = [double]"7.584431349227927e-05"
$ws.Range("AH8").Value = [double]"-1.190835220626363e-07"
$ws.Range("AI8").Value = [double]"1.323732989355032e-06"

$ws.Range("AF9").Value = [double]"8.56496544656693e-05"
$ws.Range("AG9").Value = [double]"-3.19439680058764e-05"
$ws.Range("AH9").Value = [double]"1.494868473621394e-06"
$ws.Range("AI9").Value = [double]"-5.575274178542705e-07"

$ws.Range("AF10").Value = [double]"-4.386145651835704e-05"
$ws.Range("AG10").Value = [double]"0.0002544637782517789"
$ws.Range("AH10").Value = [double]"-7.655268309656591e-07"
$ws.Range("AI10").Value = [double]"4.441230757558282e-06"

$ws.Range("AF11").Value = [double]"-1.449225119642961e-05"
$ws.Range("AG11").Value = [double]"-8.437436580521762e-07"
$ws.Range("AH11").Value = [double]"-2.529374994037841e-07"
$ws.Range("AI11").Value = [double]"-1.472610487583164e-08"

$ws.Range("AF12").Value = [double]"-7.470480275628688e-05"
$ws.Range("AG12").Value = [double]"5.878220435562298e-05"
$ws.Range("AH12").Value = [double]"-1.303844775150141e-06"
$ws.Range("AI12").Value = [double]"1.025943007585773e-06"

$ws.Range("AF13").Value = [double]"2.975008300154514e-05"
$ws.Range("AG13").Value = [double]"1.315898747300537e-06"
$ws.Range("AH13").Value = [double]"5.192369011185599e-07"
$ws.Range("AI13").Value = [double]"2.296676576326321e-08"

$ws.Range("AF14").Value = [double]"-6.342351402111035e-05"
$ws.Range("AG14").Value = [double]"-3.033503335725385e-05"
$ws.Range("AH14").Value = [double]"-1.106949142853164e-06"
$ws.Range("AI14").Value = [double]"-5.29446210786389e-07"

$ws.Range("AF15").Value = [double]"-0.0001020770944109017"
$ws.Range("AG15").Value = [double]"-5.43551590688196e-05"
$ws.Range("AH15").Value = [double]"-1.781581388339337e-06"
$ws.Range("AI15").Value = [double]"-9.486764911961571e-07"

$ws.Range("AF16").Value = [double]"-1.621439918153555e-05"
$ws.Range("AG16").Value = [double]"1.766332894037603e-05"
$ws.Range("AH16").Value = [double]"-2.829946519504692e-07"
$ws.Range("AI16").Value = [double]"3.08283246872363e-07"

$ws.Range("AF17").Value = [double]"-6.163220461985475e-05"
$ws.Range("AG17").Value = [double]"6.955269239661277e-06"
$ws.Range("AH17").Value = [double]"-1.075684895879325e-06"
$ws.Range("AI17").Value = [double]"1.213923485947718e-07"

$ws.Range("AF18").Value = [double]"1.703056024382477e-05"
$ws.Range("AG18").Value = [double]"-0.0001763776441561227"
$ws.Range("AH18").Value = [double]"2.972393497139904e-07"
$ws.Range("AI18").Value = [double]"-3.078370617435276e-06"

$ws.Range("AF19").Value = [double]"-1.028692199156467e-05"
$ws.Range("AG19").Value = [double]"-6.264204433747977e-05"
$ws.Range("AH19").Value = [double]"-1.795406586486158e-07"
$ws.Range("AI19").Value = [double]"-1.093309923869292e-06"

$ws.Range("AF20").Value = [double]"-4.43029632890557e-05"
$ws.Range("AG20").Value = [double]"-6.752761438466592e-05"
$ws.Range("AH20").Value = [double]"-7.732325777841981e-07"
$ws.Range("AI20").Value = [double]"-1.178579207029505e-06"

$ws.Range("AF21").Value = [double]"-7.344899216832346e-06"
$ws.Range("AG21").Value = [double]"-1.79581201118495e-05"
$ws.Range("AH21").Value = [double]"-1.281926745608773e-07"
$ws.Range("AI21").Value = [double]"-3.134283234203861e-07"

$ws.Range("AF22").Value = [double]"-6.752487379912964e-06"
$ws.Range("AG22").Value = [double]"4.141845690064372e-05"
$ws.Range("AH22").Value = [double]"-1.178531374788464e-07"
$ws.Range("AI22").Value = [double]"7.228884440115987e-07"

$ws.Range("AF23").Value = [double]"2.324642542816946e-05"
$ws.Range("AG23").Value = [double]"-6.342008811088817e-05"
$ws.Range("AH23").Value = [double]"4.057266630408897e-07"
$ws.Range("AI23").Value = [double]"-1.106889349439909e-06"

$ws.Range("AF24").Value = [double]"3.448616863011011e-05"
$ws.Range("AG24").Value = [double]"3.668124846978316e-05"
$ws.Range("AH24").Value = [double]"6.018971889934039e-07"
$ws.Range("AI24").Value = [double]"6.402085595398478e-07"
